$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trading log rows appended at the end (rows 168-171)
$rows = @(
    @{ Row = 168; A = "2026-01-07T01:48:47.401656"; B = "TRADING_ATTEMPT"; C = "ETH"; D = "UNKNOWN"; E = 3246.247548899048; K = "ATTEMPT"; L = "Attempting trade 1/2" },
    @{ Row = 169; A = "2026-01-07T01:48:48.952047"; B = "POSITION_FAILED"; C = "ETH"; D = "UNKNOWN"; E = $null;              K = "FAILED";  L = "Trade execution failed for trade 1" },
    @{ Row = 170; A = "2026-01-07T01:48:49.006777"; B = "TRADING_ATTEMPT"; C = "SUI"; D = "UNKNOWN"; E = 1.861190696770879; K = "ATTEMPT"; L = "Attempting trade 2/2" },
    @{ Row = 171; A = "2026-01-07T01:48:51.004658"; B = "POSITION_FAILED"; C = "SUI"; D = "UNKNOWN"; E = $null;              K = "FAILED";  L = "Trade execution failed for trade 2" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    if ($null -ne $r.E) {
        $ws.Cells.Item($rowIndex, 5).Value = $r.E
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = ""
    }
    $ws.Cells.Item($rowIndex, 6).Value = ""
    $ws.Cells.Item($rowIndex, 7).Value = ""
    $ws.Cells.Item($rowIndex, 8).Value = ""
    $ws.Cells.Item($rowIndex, 9).Value = ""
    $ws.Cells.Item($rowIndex, 10).Value = ""
    $ws.Cells.Item($rowIndex, 11).Value = $r.K
    $ws.Cells.Item($rowIndex, 12).Value = $r.L
}
